# Insert a new "04dec2025" data column (column E) into both worksheets,
# shifting the existing 25nov2025..30nov2025 columns (E:J) one column to
# the right (F:K), and populate the new column with its values.

$wb = $excel.ActiveWorkbook

# row -> value for the new "04dec2025" column, keyed by worksheet row number
$newValues = @{
    2 = 32
    3 = 37
    4 = 35
    5 = 21
    6 = 20
    7 = 19
    8 = 13
    9 = 15
    10 = 28
    11 = 15
    12 = 17
    13 = 36
}

foreach ($sheetName in @("crosstab", "annot")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert a new column at E, shifting old E:J -> F:K
    $ws.Columns("E:E").Insert()

    # Header cell - plain text, matches style of the row (inherited from Insert)
    $ws.Cells.Item(1, 5).Value = "04dec2025"

    foreach ($row in $newValues.Keys) {
        $value = $newValues[$row]
        $cell = $ws.Cells.Item($row, 5)

        if ($sheetName -eq "annot") {
            # This sheet stores its numbers as text strings; force text storage
            # then restore the default (unstyled) cell formatting so the
            # written cell matches the style of its neighboring data cells.
            $cell.NumberFormat = "@"
            $cell.Value = [string]$value
            $cell.Style = "Normal"
        }
        else {
            # This sheet stores genuine numeric values
            $cell.Value = $value
        }
    }
}
